# Refresh market-price derived figures across the Leve profit sheets.
# Generated from the authoritative cell-level diff (sheet/cell -> new value or clear).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 446
$ws.Range("I2").Value = 224.4
$ws.Range("K2").Value = 224.4
$ws.Range("M2").Value = -111.4
$ws.Range("H39").Value = 91.94118
$ws.Range("I39").Value = 91.94118
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 275.82354
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 20.17646000000002
$ws.Range("N39").ClearContents()
$ws.Range("H40").Value = 1971.75
$ws.Range("J40").Value = 2096.3333
$ws.Range("L40").Value = 2096.3333
$ws.Range("N40").Value = -2446.3333
$ws.Range("H98").Value = 842.46155
$ws.Range("I98").Value = 860.16
$ws.Range("K98").Value = 860.16
$ws.Range("M98").Value = 637.84
$ws.Range("H112").Value = 2903.4736
$ws.Range("I112").Value = 1831.6666
$ws.Range("K112").Value = 5494.9998
$ws.Range("M112").Value = -4386.9998
$ws.Range("H115").Value = 765.8333
$ws.Range("I115").Value = 765.8333
$ws.Range("K115").Value = 2297.4999
$ws.Range("M115").Value = -730.4998999999998
$ws.Range("H122").Value = 842.46155
$ws.Range("I122").Value = 860.16
$ws.Range("K122").Value = 2580.48
$ws.Range("M122").Value = -130.48
$ws.Range("H137").Value = 8608.666999999999
$ws.Range("I137").Value = 2817.2856
$ws.Range("J137").Value = 13676.125
$ws.Range("K137").Value = 8451.856800000001
$ws.Range("L137").Value = 41028.375
$ws.Range("M137").Value = -5901.856800000001
$ws.Range("N137").Value = -46128.375
$ws.Range("H138").Value = 5128.857
$ws.Range("J138").Value = 4488.857
$ws.Range("L138").Value = 13466.571
$ws.Range("N138").Value = -23746.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 534.9231
$ws.Range("J5").Value = 793.5
$ws.Range("L5").Value = 793.5
$ws.Range("N5").Value = -1017.5
$ws.Range("H32").Value = 198333.67
$ws.Range("I32").Value = 212342.06
$ws.Range("K32").Value = 212342.06
$ws.Range("M32").Value = -212055.06
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H45").Value = 6000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 6000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -6754
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H74").Value = 13936.944
$ws.Range("I74").Value = 7912.5
$ws.Range("J74").Value = 15658.214
$ws.Range("K74").Value = 7912.5
$ws.Range("L74").Value = 15658.214
$ws.Range("M74").Value = -7038.5
$ws.Range("N74").Value = -17406.214
$ws.Range("H77").Value = 13936.944
$ws.Range("I77").Value = 7912.5
$ws.Range("J77").Value = 15658.214
$ws.Range("K77").Value = 39562.5
$ws.Range("L77").Value = 78291.07000000001
$ws.Range("M77").Value = -35194.5
$ws.Range("N77").Value = -87027.07000000001
$ws.Range("H102").Value = 1767.1
$ws.Range("I102").Value = 1771.375
$ws.Range("J102").Value = 1750
$ws.Range("K102").Value = 1771.375
$ws.Range("L102").Value = 1750
$ws.Range("M102").Value = -149.375
$ws.Range("N102").Value = -4994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 534.9231
$ws.Range("J4").Value = 793.5
$ws.Range("L4").Value = 793.5
$ws.Range("N4").Value = -1023.5
$ws.Range("H22").Value = 799
$ws.Range("I22").Value = 498.33334
$ws.Range("K22").Value = 498.33334
$ws.Range("M22").Value = -325.33334
$ws.Range("H86").Value = 3109.6924
$ws.Range("I86").Value = 2742.6
$ws.Range("K86").Value = 2742.6
$ws.Range("M86").Value = -1619.6
$ws.Range("H89").Value = 3109.6924
$ws.Range("I89").Value = 2742.6
$ws.Range("K89").Value = 13713
$ws.Range("M89").Value = -8097
$ws.Range("H94").Value = 1877.129
$ws.Range("I94").Value = 1310.0358
$ws.Range("J94").Value = 7170
$ws.Range("K94").Value = 1310.0358
$ws.Range("L94").Value = 7170
$ws.Range("M94").Value = -859.0358000000001
$ws.Range("N94").Value = -8072
$ws.Range("H105").Value = 2203.2812
$ws.Range("I105").Value = 1939
$ws.Range("K105").Value = 1939
$ws.Range("M105").Value = -192

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 391441.16
$ws.Range("I22").Value = 607627.3
$ws.Range("J22").Value = 2306
$ws.Range("K22").Value = 607627.3
$ws.Range("L22").Value = 2306
$ws.Range("M22").Value = -607277.3
$ws.Range("N22").Value = -3006
$ws.Range("H31").Value = 1769.32
$ws.Range("I31").Value = 1749.6316
$ws.Range("J31").Value = 1831.6666
$ws.Range("K31").Value = 1749.6316
$ws.Range("L31").Value = 1831.6666
$ws.Range("M31").Value = -1454.6316
$ws.Range("N31").Value = -2421.6666
$ws.Range("H34").Value = 1769.32
$ws.Range("I34").Value = 1749.6316
$ws.Range("J34").Value = 1831.6666
$ws.Range("K34").Value = 1749.6316
$ws.Range("L34").Value = 1831.6666
$ws.Range("M34").Value = -1547.6316
$ws.Range("N34").Value = -2235.6666
$ws.Range("H60").Value = 26299.8
$ws.Range("I60").Value = 23500
$ws.Range("J60").Value = 28166.334
$ws.Range("K60").Value = 23500
$ws.Range("L60").Value = 28166.334
$ws.Range("M60").Value = -22989
$ws.Range("N60").Value = -29188.334
$ws.Range("H99").Value = 13784
$ws.Range("J99").Value = 3124.3
$ws.Range("L99").Value = 3124.3
$ws.Range("N99").Value = -6120.3
$ws.Range("H126").Value = 13784
$ws.Range("J126").Value = 3124.3
$ws.Range("L126").Value = 9372.900000000001
$ws.Range("N126").Value = -14312.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4000087.5
$ws.Range("I2").Value = 5263198.5
$ws.Range("K2").Value = 5263198.5
$ws.Range("M2").Value = -5263085.5
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H97").Value = 1207.2307
$ws.Range("I97").Value = 1233.3
$ws.Range("K97").Value = 1233.3
$ws.Range("M97").Value = -737.3
$ws.Range("H126").Value = 3322.375
$ws.Range("I126").Value = 3090.6667
$ws.Range("K126").Value = 9272.000100000001
$ws.Range("M126").Value = -6802.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9309
$ws.Range("H136").Value = 3429.8635
$ws.Range("I136").Value = 3087.4443
$ws.Range("J136").Value = 4970.75
$ws.Range("K136").Value = 9262.332900000001
$ws.Range("L136").Value = 14912.25
$ws.Range("M136").Value = -6712.332900000001
$ws.Range("N136").Value = -20012.25
